$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("Nodes-Number of edg.").Name = "Nodes-Number of edge."
$wb.Worksheets.Item("Nodes-Number of edg. (Fit Data)").Name = "Nodes-Number of edge. (Fitting)"
$wb.Worksheets.Item("Nodes-Average degree (Fit Data)").Name = "Nodes-Average degree (Fitting)"
$wb.Worksheets.Item("Nodes-Graph density (Fit Data)").Name = "Nodes-Graph density (Fitting)"
